# Daily attendance processing - 2025-10-07 09:20:06
#
# - Reorders the "Recorded By" (column G) attributions on a batch of rows.
# - Refreshes the derived statistics: the "Class Statistics" box (K/L) and
#   the "Group Statistics" table (K:S) for the B2D/B2E/B2F rows.
# - Two sessions (rows 95, 121, 147) moved from "Pending" (0 students, no
#   recorder) to "Recorded" (System recorded a first batch of students),
#   which also flips their row styling from the pending/yellow look to the
#   recorded/green look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# "Recorded By" (column G) - reorder "system, backup@backdoor.com, System"
# -> "backup@backdoor.com, System, system"
# ---------------------------------------------------------------------
$backdoorRows = @(2, 29, 56)
foreach ($r in $backdoorRows) {
    $ws.Cells.Item($r, 7).Value = "backup@backdoor.com, System, system"
}

# ---------------------------------------------------------------------
# "Recorded By" (column G) - reorder "System, dnasr281@gmail.com"
# -> "dnasr281@gmail.com, System"
# ---------------------------------------------------------------------
$dnasrRows = @(3, 6, 11, 12, 13, 14, 15, 30, 33, 38, 39, 40, 41, 42, 57, 60, 65, 66, 67, 68, 69, 89, 93, 115, 119, 141, 145)
foreach ($r in $dnasrRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------
# Class Statistics box (K/L columns, B2A block)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 83          # Recorded Sessions
$ws.Range("L8").Value = 75          # Pending Sessions
$ws.Range("L9").Value = "'52.2%"    # Coverage %  (leading ' keeps it literal text)
$ws.Range("L10").Value = "'67.1%"   # Average Attendance %

# ---------------------------------------------------------------------
# Group Statistics table (K:S columns) - rows for B2D, B2E, B2F
# ---------------------------------------------------------------------
$ws.Range("O18").Value = 13
$ws.Range("Q18").Value = 13
$ws.Range("R18").Value = "'50.0%"
$ws.Range("S18").Value = "'69.0%"

$ws.Range("O19").Value = 13
$ws.Range("Q19").Value = 13
$ws.Range("R19").Value = "'50.0%"
$ws.Range("S19").Value = "'70.3%"

$ws.Range("O20").Value = 13
$ws.Range("Q20").Value = 13
$ws.Range("R20").Value = "'50.0%"
$ws.Range("S20").Value = "'67.5%"

# ---------------------------------------------------------------------
# Sessions newly recorded by "System" (rows 95, 121, 147) - were "Pending"
# with 0 students; now have a first batch of attendance recorded. Copy the
# "Recorded" row formatting from the row immediately above (already styled
# green) onto these rows, then fill in Recorded By / Students / Status.
# ---------------------------------------------------------------------
$ws.Range("A94:I94").Copy()
$ws.Range("A95:I95").PasteSpecial(-4122)
$ws.Cells.Item(95, 7).Value = "System"
$ws.Cells.Item(95, 8).Value = "1/56"
$ws.Cells.Item(95, 9).Value = "Recorded"

$ws.Range("A120:I120").Copy()
$ws.Range("A121:I121").PasteSpecial(-4122)
$ws.Cells.Item(121, 7).Value = "System"
$ws.Cells.Item(121, 8).Value = "2/55"
$ws.Cells.Item(121, 9).Value = "Recorded"

$ws.Range("A146:I146").Copy()
$ws.Range("A147:I147").PasteSpecial(-4122)
$ws.Cells.Item(147, 7).Value = "System"
$ws.Cells.Item(147, 8).Value = "2/57"
$ws.Cells.Item(147, 9).Value = "Recorded"

$excel.CutCopyMode = 0
